# Updated cryptos list on Mon May 15 12:30:05 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures on
# the active sheet for the latest crypto snapshot. Every cell in D/E holds
# plain text (not a real number/percentage) in the source data, so values
# are written as text and, for column D, the cell is briefly switched to a
# text number-format so Excel doesn't silently re-interpret price strings
# such as "0.4315" or "21.61" as floating point numbers; the style is then
# restored to "Normal" so no formatting change is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-PriceText "D2" "27.689.04"
$ws.Range("E2").Value = "  -0.28%  "

Set-PriceText "D3" "1.848.95"
$ws.Range("E3").Value = "  -0.83%  "

$ws.Range("E4").Value = "  -2.72%  "

Set-PriceText "D5" "320.38"
$ws.Range("E5").Value = "  -1.22%  "

$ws.Range("E6").Value = "  -2.44%  "

Set-PriceText "D7" "0.4315"
$ws.Range("E7").Value = "  -2.55%  "

Set-PriceText "D8" "0.3745"
$ws.Range("E8").Value = "  -1.38%  "

Set-PriceText "D9" "0.07353"
$ws.Range("E9").Value = "  -1.55%  "

Set-PriceText "D10" "0.8802"
$ws.Range("E10").Value = "  -0.44%  "

Set-PriceText "D11" "21.61"
$ws.Range("E11").Value = "  -0.10%  "

Set-PriceText "D12" "1.863.64"
$ws.Range("E12").Value = "  -0.42%  "

$ws.Range("E13").Value = "  -0.50%  "

Set-PriceText "D14" "5.455"
$ws.Range("E14").Value = "  -2.04%  "

Set-PriceText "D15" "0.07149"
$ws.Range("E15").Value = "  -1.28%  "

Set-PriceText "D16" "87.77"
$ws.Range("E16").Value = "  +4.74%  "

Set-PriceText "D17" "1.016"
$ws.Range("E17").Value = "  -2.62%  "

Set-PriceText "D18" "0.000009004"
$ws.Range("E18").Value = "  -1.49%  "

Set-PriceText "D19" "1.013"
$ws.Range("E19").Value = "  -2.39%  "

Set-PriceText "D20" "15.47"
$ws.Range("E20").Value = "  -0.61%  "

Set-PriceText "D21" "27.701.49"
$ws.Range("E21").Value = "  -0.28%  "

Set-PriceText "D22" "5.246"
$ws.Range("E22").Value = "  -1.45%  "

$ws.Range("E23").Value = "  -1.56%  "

Set-PriceText "D24" "2.089.43"
$ws.Range("E24").Value = "  -0.60%  "

Set-PriceText "D25" "2.008"
$ws.Range("E25").Value = "  -0.03%  "

Set-PriceText "D26" "155.86"
$ws.Range("E26").Value = "  -1.96%  "

$ws.Range("E27").Value = "  -1.23%  "

Set-PriceText "D28" "2.109"
$ws.Range("E28").Value = "  +6.51%  "

Set-PriceText "D29" "5.393"
$ws.Range("E29").Value = "  +1.40%  "

Set-PriceText "D30" "120.09"
$ws.Range("E30").Value = "  +1.81%  "

Set-PriceText "D31" "0.08944"
$ws.Range("E31").Value = "  -1.31%  "

Set-PriceText "D32" "1.236"
$ws.Range("E32").Value = "  +1.91%  "

Set-PriceText "D33" "0.7764"

Set-PriceText "D34" "4.565"
$ws.Range("E34").Value = "  +0.03%  "

Set-PriceText "D35" "2.924"
$ws.Range("E35").Value = "  -3.64%  "

$ws.Range("E36").Value = "  -2.53%  "

Set-PriceText "D37" "1.137"
$ws.Range("E37").Value = "  -1.23%  "

Set-PriceText "D38" "0.05340"
$ws.Range("E38").Value = "  +0.05%  "

$ws.Range("E39").Value = "  -0.84%  "

Set-PriceText "D40" "7.186"
$ws.Range("E40").Value = "  +4.71%  "

Set-PriceText "D41" "2.872"
$ws.Range("E41").Value = "  +0.09%  "

Set-PriceText "D42" "0.5159"
$ws.Range("E42").Value = "  -0.70%  "

$ws.Range("E43").Value = "  -0.66%  "

Set-PriceText "D44" "8.823"
$ws.Range("E44").Value = "  +1.97%  "

Set-PriceText "D45" "109.28"
$ws.Range("E45").Value = "  -0.85%  "

Set-PriceText "D46" "10.65"
$ws.Range("E46").Value = "  -0.15%  "

Set-PriceText "D47" "0.4726"
$ws.Range("E47").Value = "  +0.46%  "

Set-PriceText "D48" "0.06502"
$ws.Range("E48").Value = "  +0.24%  "

Set-PriceText "D49" "1.698"
$ws.Range("E49").Value = "  -1.19%  "

Set-PriceText "D50" "1.014"
$ws.Range("E50").Value = "  -2.64%  "

Set-PriceText "D51" "1.872"
$ws.Range("E51").Value = "  -1.95%  "
